# Insert a blank paragraph, a "Grammar Test Paragraph:" heading, and a
# paragraph containing several grammar errors right after the document's
# opening title paragraph ("Test Document for Document Processor").

$d = $word.ActiveDocument

# The title is the first paragraph in the document.
$titlePara = $d.Paragraphs.Item(1)

# InsertParagraphAfter() splits in a new, empty paragraph right after the
# anchor paragraph's range - calling it three times (re-fetching the
# paragraph collection each time) lays down three fresh empty paragraphs
# immediately following the title.
$titlePara.Range.InsertParagraphAfter()
$d.Paragraphs.Item(2).Range.InsertParagraphAfter()
$d.Paragraphs.Item(3).Range.InsertParagraphAfter()

# Paragraph 2 stays blank (matches the spacer paragraph already used
# elsewhere in the document) - set it explicitly so it gets a proper
# empty <w:t/> run instead of a run with no text element at all.
$d.Paragraphs.Item(2).Range.Text = ""

# Paragraph 3: the new section heading.
$d.Paragraphs.Item(3).Range.Text = "Grammar Test Paragraph:"

# Paragraph 4: the sentence riddled with grammar/spelling errors used to
# exercise the grammar checker in the document processor tests.
$d.Paragraphs.Item(4).Range.Text = "The cats and dog is running fast. We dont need no help with grammer. This sentense contains muliple mispelled words. The weather have been nice yesterday?"
